# Daily attendance processing - swap "System, <email>" to "<email>, System"
# in column G (Edited By) across all data rows of the session analysis sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1

$count = 0
for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $text = $cell.Text
    if ($text -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
        $count = $count + 1
    }
}

Write-Host "Updated $count cells in column G"
